$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1 - "Final project check-in"
# Merge the two runs of "Post links to your procedures or " / "printed
# materials" back into a single run (content placeholder, paragraph 2).
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
$para1_2 = $tr1.Paragraphs(2, 1)
# Route through an unrelated placeholder string first so the text-diff
# engine behind TextRange.Text rewrites the paragraph as one clean run
# instead of patching a shared prefix/suffix into two runs.
$para1_2.Text = "placeholder"
$para1_2.Text = "Post links to your procedures or printed materials"

# ---------------------------------------------------------------
# Slide 12 - "Final project presentations"
# "Aim for 5-8 minutes" -> "Aim for 10-12 minutes"
# "5 groups/hour with questions" -> "3 " + "groups/hour with questions"
# ---------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)
$tr12 = $sh12.TextFrame.TextRange

$para12_9 = $tr12.Paragraphs(9, 1)
$para12_9.Text = "placeholder"
$para12_9.Text = "Aim for 10-12 minutes"

$para12_10 = $tr12.Paragraphs(10, 1)
$lead = $para12_10.Characters(1, 2)
$lead.Text = "3 "

# ---------------------------------------------------------------
# Slide 6 - "Case Study Research Methods"
# "Non-experimental" -> "Observational, non-experimental"
# "Drive theory and subsequent experimental research" ->
#   "Generate theory and inspire subsequent experimental research"
# ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange

$para6_2 = $tr6.Paragraphs(2, 1)
$para6_2.Text = "placeholder"
$para6_2.Text = "Observational, non-experimental"

$para6_4 = $tr6.Paragraphs(4, 1)
$para6_4.Text = "placeholder"
$para6_4.Text = "Generate theory and inspire subsequent experimental research"
